$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the worksheet; Excel automatically updates defined names and
# chart series formulas (c:f) that reference the sheet by name.
$ws.Name = "Operations"

# Update the selection on the sheet to match the new active cell.
$ws.Range("C41").Select()
